$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2..23) down by one row (to 3..24), working
# from the bottom up so we don't overwrite data before it is copied.
for ($r = 23; $r -ge 2; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $dst.Value = $src.Value2
}

# The row shifted into row 24 is brand new on the sheet, so it has no
# formatting yet. Copy the date-column number format from row 23 (the
# row the data originally came from) so it renders the same as before.
$ws.Range("D24").NumberFormat = $ws.Range("D23").NumberFormat

# Insert the new weekly record into row 2 (newest date first).
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44515
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112026
$ws.Range("G2").Value = "Haba"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7500
$ws.Range("N2").Value = '$/saco 25 kilos'
$ws.Range("O2").Value = "Provincia de Diguillín"
$ws.Range("P2").Value = 300
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"

Write-Output "done"
